$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 34 was previously styled as the "last row" (date-only format).
# Now that a new row is being appended, row 34 becomes a regular row,
# so copy the number format used by the other regular rows (e.g. A33).
$ws.Range("A34").NumberFormat = $ws.Range("A33").NumberFormat

# Append the new day's data as row 35, using the style previously used
# by row 34 (the "last row" style, date-only format).
$ws.Range("A35").NumberFormat = "YYYY-MM-DD"

$ws.Range("A35").Value = 45619
$ws.Range("B35").Value = 89
$ws.Range("C35").Value = 73
$ws.Range("D35").Value = 82
